# Atualização de bases das ligas, do dia: 19-06-2024 às 21:51
#
# The upstream data refresh reshuffled a handful of already-scraped match
# rows (same "Thailand Premier League" table, same id set) into a
# different row order. Concretely:
#   - rows 117 / 118 swap their entire record (columns B..AD, i.e. every
#     column except the running index in column A)
#   - rows 234 / 235 / 236 cyclically rotate their entire record the same
#     way (234<-235, 235<-236, 236<-234)
#
# We implement this generically: read each row's B..AD values off the
# live sheet first (so we are moving whatever is actually there rather
# than hard-coding the destination), then write them back in the new
# order. Column A (the sequential row index) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns present on every data row.
$allCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

function Read-RowValues($row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range($c + $row).Value()
    }
    return $vals
}

function Write-RowValues($row, $cols, $vals) {
    foreach ($c in $cols) {
        $v = $vals[$c]
        if ($v -ne $null) {
            $ws.Range($c + $row).Value = $v
        }
    }
}

# --- Swap rows 117 and 118 (full record, columns B..AD) -------------------
$row117 = Read-RowValues 117 $allCols
$row118 = Read-RowValues 118 $allCols

Write-RowValues 117 $allCols $row118
Write-RowValues 118 $allCols $row117

# --- Cyclic rotation of rows 234, 235, 236 (full record, columns B..AD) ---
# Note: rows 234-236 don't populate I/J (HTHG/HTAG not yet known), so we
# only touch the columns that are actually used on these rows.
$rotCols = @("B","C","D","E","F","G","H","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$row234 = Read-RowValues 234 $rotCols
$row235 = Read-RowValues 235 $rotCols
$row236 = Read-RowValues 236 $rotCols

# 234 <- 235, 235 <- 236, 236 <- 234
Write-RowValues 234 $rotCols $row235
Write-RowValues 235 $rotCols $row236
Write-RowValues 236 $rotCols $row234
